$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.428.58'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.624.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.06%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.98'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.83'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +9.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.624.67'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000188'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.53%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.335.34'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.08%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.623.68'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.02'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '380.54'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.60'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.19'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.56%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +18.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.20'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.84%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.38'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.97'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +9.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.761.36'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.29%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0956'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '520.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.08'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.58%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.75%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '165.44'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.33'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.41'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.51%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.08'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -7.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.07'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.40%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +8.61%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.03%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.75'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.543'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.00%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +7.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0264'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.58%  '
